$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 columns A-H become text-typed cells (stored as strings even though
# their content looks numeric) as part of new input-validation/error
# handling for the Model Settings row. A gets "12", H gets "12" (was 100).
$ws.Range("A2:H2").NumberFormat = "@"

$ws.Range("A2").Value = "12"
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "12"

# I2 (Model) was an empty inline string; it now holds numeric 0.
$ws.Range("I2").Value = 0
